$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 224
$ws1.Range("F4").Value = 762
$ws1.Range("F6").Value = 397
$ws1.Range("F7").Value = 544
$ws1.Range("F11").Value = 122
$ws1.Range("F12").Value = 585
$ws1.Range("F13").Value = 76
$ws1.Range("F14").Value = 1741
$ws1.Range("F15").Value = 310
$ws1.Range("F16").Value = 2005
$ws1.Range("F17").Value = 242
$ws1.Range("F19").Value = 37
$ws1.Range("F20").Value = 128

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 219
$ws2.Range("F5").Value = 11
$ws2.Range("F11").Value = 13
$ws2.Range("F13").Value = 82

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5265
$ws3.Range("F3").Value = 300
$ws3.Range("F4").Value = 121

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5265
$ws4.Range("F4").Value = 300
$ws4.Range("F6").Value = 121
$ws4.Range("F7").Value = 224
$ws4.Range("F8").Value = 219
$ws4.Range("F10").Value = 11
$ws4.Range("F13").Value = 763
$ws4.Range("F17").Value = 397
$ws4.Range("F18").Value = 544
$ws4.Range("F23").Value = 122
$ws4.Range("F24").Value = 13
$ws4.Range("F26").Value = 585
$ws4.Range("F27").Value = 76
$ws4.Range("F28").Value = 82
$ws4.Range("F29").Value = 1741
$ws4.Range("F30").Value = 310
$ws4.Range("F31").Value = 2005
$ws4.Range("F33").Value = 242
$ws4.Range("F35").Value = 37
$ws4.Range("F36").Value = 128
